$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Column C ("Förändrad") holds a date serial for every data row (rows 2-387).
# The workbook was refreshed a day later, so bump each date by one day
# (45177 -> 45178), leaving every other cell untouched.
$rng = $ws.Range("C2:C387")
$vals = $rng.Value2
$rows = $vals.GetLength(0)
for ($i = 1; $i -le $rows; $i++) {
    if ($vals[$i, 1] -eq 45177) {
        $vals[$i, 1] = 45178
    }
}
$rng.Value2 = $vals
